# Update the "dSF" column (F) values for specific rows to repull/push data
# and reflect the corrected mean calculation, per commit message:
# "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = 1
$ws.Range("F7").Value = -1
$ws.Range("F9").Value = 2
$ws.Range("F17").Value = 2
$ws.Range("F18").Value = 3
$ws.Range("F20").Value = 1
$ws.Range("F26").Value = 1
$ws.Range("F27").Value = -1
$ws.Range("F30").Value = 6
$ws.Range("F32").Value = 2
$ws.Range("F33").Value = -1
$ws.Range("F38").Value = 2
$ws.Range("F40").Value = -2
$ws.Range("F41").Value = 2
$ws.Range("F42").Value = 1
$ws.Range("F43").Value = 0
$ws.Range("F49").Value = 1
$ws.Range("F50").Value = -4
$ws.Range("F56").Value = -4
$ws.Range("F59").Value = 0
$ws.Range("F63").Value = -2
$ws.Range("F64").Value = 0
$ws.Range("F65").Value = -2
$ws.Range("F69").Value = -1
$ws.Range("F72").Value = 0
$ws.Range("F73").Value = -1
$ws.Range("F76").Value = -1
